# Auto-generated Excel COM-interop script to update the cryptos list
# (values refreshed by the scheduled GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2:E51 holds price/volume text that looks numeric (e.g. '0.4765', '29.512.06').
# Force the range to Text format first so Excel stores the exact original
# string instead of silently coercing it to a floating point Double, then
# clear the temporary formatting again so no stray cell styles are left behind.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.512.06'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '1.899.51'
$ws.Range("E3").Value = '  -0.92%  '
$ws.Range("E4").Value = '  -0.60%  '
$ws.Range("D5").Value = '339.00'
$ws.Range("E5").Value = '  +4.23%  '
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").Value = '0.4765'
$ws.Range("E7").Value = '  -1.07%  '
$ws.Range("D8").Value = '0.4011'
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("D9").Value = '47.22'
$ws.Range("E9").Value = '  -1.77%  '
$ws.Range("D10").Value = '0.08039'
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("D11").Value = '0.9926'
$ws.Range("E11").Value = '  -1.75%  '
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").Value = '1.908.60'
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").Value = '5.943'
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("D15").Value = '7.093'
$ws.Range("E15").Value = '  -1.83%  '
$ws.Range("D16").Value = '89.18'
$ws.Range("E16").Value = '  -2.61%  '
$ws.Range("D17").Value = '0.06794'
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").Value = '0.00001021'
$ws.Range("E19").Value = '  -1.77%  '
$ws.Range("E20").Value = '  -1.29%  '
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").Value = '29.508.76'
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("E23").Value = '  -2.54%  '
$ws.Range("D24").Value = '11.67'
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("D25").Value = '2.153'
$ws.Range("E25").Value = '  -1.63%  '
$ws.Range("D26").Value = '2.155.74'
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").Value = '157.50'
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("D28").Value = '6.520'
$ws.Range("E28").Value = '  -0.99%  '
$ws.Range("D29").Value = '19.67'
$ws.Range("E29").Value = '  -1.55%  '
$ws.Range("E30").Value = '  -2.60%  '
$ws.Range("D31").Value = '119.23'
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("D32").Value = '0.9983'
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("D33").Value = '0.09559'
$ws.Range("E33").Value = '  -0.64%  '
$ws.Range("D34").Value = '5.488'
$ws.Range("E34").Value = '  -2.68%  '
$ws.Range("D35").Value = '1.389'
$ws.Range("E35").Value = '  +1.24%  '
$ws.Range("D36").Value = '3.531'
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").Value = '0.06425'
$ws.Range("E37").Value = '  +5.37%  '
$ws.Range("D38").Value = '0.02246'
$ws.Range("E38").Value = '  -1.61%  '
$ws.Range("D39").Value = '1.201'
$ws.Range("E39").Value = '  +1.41%  '
$ws.Range("D40").Value = '0.5841'
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("D41").Value = '10.57'
$ws.Range("E41").Value = '  -2.84%  '
$ws.Range("D42").Value = '7.762'
$ws.Range("E42").Value = '  -3.60%  '
$ws.Range("D43").Value = '0.1823'
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("D44").Value = '2.430'
$ws.Range("E44").Value = '  +2.03%  '
$ws.Range("D45").Value = '1.265'
$ws.Range("E45").Value = '  -1.14%  '
$ws.Range("D46").Value = '12.18'
$ws.Range("E46").Value = '  -1.61%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.5508'
$ws.Range("E47").Value = '  -1.35%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.07338'
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("D49").Value = '1.957'
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("D50").Value = '116.57'
$ws.Range("E50").Value = '  -1.62%  '
$ws.Range("D51").Value = '71.24'
$ws.Range("E51").Value = '  -1.23%  '

# Restore the original (default) cell formatting on the range we touched.
$priceRange.ClearFormats()

